$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 196, pushing existing rows 196-288 down to 198-290.
$ws.Range("A196:A197").EntireRow.Insert()

# --- New row 196 ---
$ws.Cells.Item(196,1).Value2 = 1
$ws.Cells.Item(196,2).Value2 = 'Agrícola del Norte S.A. de Arica'
$ws.Cells.Item(196,3).Value2 = 'Arica y Parinacota'
$ws.Cells.Item(196,4).Value2 = 44839
$ws.Cells.Item(196,5).Value2 = 15
$ws.Cells.Item(196,6).Value2 = 'Fruta'
$ws.Cells.Item(196,7).Value2 = 100102
$ws.Cells.Item(196,8).Value2 = 'Cítricos'
$ws.Cells.Item(196,9).Value2 = 100102003
$ws.Cells.Item(196,10).Value2 = 'Limón'
$ws.Cells.Item(196,11).Value2 = 'Sin especificar'
$ws.Cells.Item(196,12).Value2 = '2a amarillo'
$ws.Cells.Item(196,13).Value2 = 300
$ws.Cells.Item(196,14).Value2 = 11000
$ws.Cells.Item(196,15).Value2 = 12000
$ws.Cells.Item(196,16).Value2 = 11500
$ws.Cells.Item(196,17).Value2 = '$/caja 20 kilos'
$ws.Cells.Item(196,18).Value2 = 'Región de Coquimbo'
$ws.Cells.Item(196,19).Value2 = 575
$ws.Cells.Item(196,20).Value2 = 20

# --- New row 197 ---
$ws.Cells.Item(197,1).Value2 = 1
$ws.Cells.Item(197,2).Value2 = 'Agrícola del Norte S.A. de Arica'
$ws.Cells.Item(197,3).Value2 = 'Arica y Parinacota'
$ws.Cells.Item(197,4).Value2 = 44839
$ws.Cells.Item(197,5).Value2 = 15
$ws.Cells.Item(197,6).Value2 = 'Fruta'
$ws.Cells.Item(197,7).Value2 = 100102
$ws.Cells.Item(197,8).Value2 = 'Cítricos'
$ws.Cells.Item(197,9).Value2 = 100102003
$ws.Cells.Item(197,10).Value2 = 'Limón'
$ws.Cells.Item(197,11).Value2 = 'Sutil De Gase'
$ws.Cells.Item(197,12).Value2 = 'Primera'
$ws.Cells.Item(197,13).Value2 = 250
$ws.Cells.Item(197,14).Value2 = 34000
$ws.Cells.Item(197,15).Value2 = 35000
$ws.Cells.Item(197,16).Value2 = 34500
$ws.Cells.Item(197,17).Value2 = '$/caja 24 kilos'
$ws.Cells.Item(197,18).Value2 = 'Perú'
$ws.Cells.Item(197,19).Value2 = 1438
$ws.Cells.Item(197,20).Value2 = 24

# Make sure the date cells keep the same number format as the rest of column D.
$ws.Cells.Item(196,4).NumberFormat = $ws.Cells.Item(198,4).NumberFormat
$ws.Cells.Item(197,4).NumberFormat = $ws.Cells.Item(198,4).NumberFormat
